# This script re-applies a row-content re-shuffle that happened in the
# source spreadsheet: whole rows (all columns A..AY) were re-assigned so
# that each destination row below ends up holding exactly the data that
# used to live in a different "source" row (values, numbers, booleans,
# text, and blanks all move together).
#
#   row 3  <- old row 4
#   row 4  <- old row 7
#   row 7  <- old row 3
#   row 10 <- old row 11
#   row 11 <- old row 10
#   row 19 <- old row 20
#   row 20 <- old row 19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 51   # column AY

# Columns that hold real numbers or booleans; every other used column in
# these rows holds text (stored as inline/shared strings in the workbook),
# even when the text looks like a number (e.g. "8") or a date
# (e.g. "2026-01-28"). Those must be written back as text, otherwise the
# COM layer "helpfully" reinterprets them as numbers/dates.
$numericCols = @(1, 2, 5, 17, 18, 19)          # A, B, E, Q, R, S
$booleanCols = @(30, 31, 33)                   # AD, AE, AG

function Is-NumericCol($c) {
    foreach ($nc in $numericCols) { if ($nc -eq $c) { return $true } }
    return $false
}
function Is-BooleanCol($c) {
    foreach ($bc in $booleanCols) { if ($bc -eq $c) { return $true } }
    return $false
}

$destToSrc = @{}
$destToSrc[3]  = 4
$destToSrc[4]  = 7
$destToSrc[7]  = 3
$destToSrc[10] = 11
$destToSrc[11] = 10
$destToSrc[19] = 20
$destToSrc[20] = 19

$rowsToSnapshot = @(3, 4, 7, 10, 11, 19, 20)

# Snapshot every involved row in full before writing anything, because
# some rows are both a source and a destination (cyclic group 3->4->7->3).
$snapshot = @{}
foreach ($r in $rowsToSnapshot) {
    $rowValues = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowValues += $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowValues
}

# Write each destination row using the snapshot captured above.
foreach ($dest in $rowsToSnapshot) {
    $src = $destToSrc[$dest]
    $data = $snapshot[$src]
    for ($c = 1; $c -le $lastCol; $c++) {
        $val = $data[$c - 1]
        $cell = $ws.Cells.Item($dest, $c)

        if (Is-NumericCol $c) {
            $cell.Value = $val
        } elseif (Is-BooleanCol $c) {
            $cell.Value = $val
        } else {
            if ([string]::IsNullOrEmpty($val)) {
                $cell.Value = $val
            } else {
                $cell.NumberFormat = "@"
                $cell.Value = [string]$val
            }
        }
    }
}
